$wb = $excel.ActiveWorkbook

# --- Rename "Requested quantity" headers on the existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet, positioned after "Monthly Trend" ---
$wsForecast = $wb.Worksheets.Add($null, $wsMonthly)
$wsForecast.Name = "PO Forecast"

# Match the page margins used by the other sheets in the workbook
# (PageSetup margins are expressed in points; 72pt = 1in)
$wsForecast.PageSetup.LeftMargin = 54
$wsForecast.PageSetup.RightMargin = 54
$wsForecast.PageSetup.TopMargin = 72
$wsForecast.PageSetup.BottomMargin = 72
$wsForecast.PageSetup.HeaderMargin = 36
$wsForecast.PageSetup.FooterMargin = 36

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Data rows
$wsForecast.Range("A2").Value = 44934.99999999999
$wsForecast.Range("A2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B2").Value = 1
$wsForecast.Range("C2").Value = 0.9999999986485325
$wsForecast.Range("D2").Value = 1.000000001252552
$wsForecast.Range("A3").Value = 44941.99999999999
$wsForecast.Range("A3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B3").Value = 1
$wsForecast.Range("C3").Value = 0.999999998827056
$wsForecast.Range("D3").Value = 1.000000001294141
$wsForecast.Range("A4").Value = 44948.99999999999
$wsForecast.Range("A4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B4").Value = 1
$wsForecast.Range("C4").Value = 0.9999999986535957
$wsForecast.Range("D4").Value = 1.000000001231421
$wsForecast.Range("A5").Value = 44962.99999999999
$wsForecast.Range("A5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B5").Value = 1
$wsForecast.Range("C5").Value = 0.9999999986930457
$wsForecast.Range("D5").Value = 1.000000001272311
$wsForecast.Range("A6").Value = 44976.99999999999
$wsForecast.Range("A6").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B6").Value = 1
$wsForecast.Range("C6").Value = 0.9999999987884297
$wsForecast.Range("D6").Value = 1.000000001289864
$wsForecast.Range("A7").Value = 44983.99999999999
$wsForecast.Range("A7").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B7").Value = 1
$wsForecast.Range("C7").Value = 0.9999999987296321
$wsForecast.Range("D7").Value = 1.000000001281928
$wsForecast.Range("A8").Value = 44990.99999999999
$wsForecast.Range("A8").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B8").Value = 1
$wsForecast.Range("C8").Value = 0.9999999985793294
$wsForecast.Range("D8").Value = 1.000000001445923
$wsForecast.Range("A9").Value = 44997.99999999999
$wsForecast.Range("A9").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B9").Value = 1
$wsForecast.Range("C9").Value = 0.9999999975207765
$wsForecast.Range("D9").Value = 1.000000002156166
$wsForecast.Range("A10").Value = 45004.99999999999
$wsForecast.Range("A10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B10").Value = 1
$wsForecast.Range("C10").Value = 0.9999999961440803
$wsForecast.Range("D10").Value = 1.000000003868573
$wsForecast.Range("A11").Value = 45011.99999999999
$wsForecast.Range("A11").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B11").Value = 1
$wsForecast.Range("C11").Value = 0.999999994175307
$wsForecast.Range("D11").Value = 1.000000006828766
$wsForecast.Range("A12").Value = 45018.99999999999
$wsForecast.Range("A12").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B12").Value = 1
$wsForecast.Range("C12").Value = 0.9999999914112131
$wsForecast.Range("D12").Value = 1.000000009566688
$wsForecast.Range("A13").Value = 45025.99999999999
$wsForecast.Range("A13").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B13").Value = 1
$wsForecast.Range("C13").Value = 0.9999999881683904
$wsForecast.Range("D13").Value = 1.000000013800566
$wsForecast.Range("A14").Value = 45032.99999999999
$wsForecast.Range("A14").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B14").Value = 1
$wsForecast.Range("C14").Value = 0.9999999848091602
$wsForecast.Range("D14").Value = 1.000000016850888
$wsForecast.Range("A15").Value = 45039.99999999999
$wsForecast.Range("A15").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B15").Value = 1
$wsForecast.Range("C15").Value = 0.999999980838043
$wsForecast.Range("D15").Value = 1.000000020650721

# Restore the originally active sheet/tab selection
$wsWeekly.Activate()

Write-Output ("Worksheets: " + $wb.Worksheets.Count)
